$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.626.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.005.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.99%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5014"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4233"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09043"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.121"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.39"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.016.14"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.078"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.482"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.012"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.32"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -6.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001115"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06681"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.963"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.632.69"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.69"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.422"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.308"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -7.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.32"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.056"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.43%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.578"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09942"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.846"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.799"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02474"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.325"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -8.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.310"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06363"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6564"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.71"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2054"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.011"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6359"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.42"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.200"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.304"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.505"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000333"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06993"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.128"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.23%  "
